# Update the "as_of_utc" timestamp column (AA) for rows 2-26
# on both the "Главные" and "Линейные" worksheets.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-10-30 07:03:06"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = $newTimestamp
    }
}
